$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028520454038103
$ws.Cells.Item(2, 4).Value = 1.036315002540566
$ws.Cells.Item(2, 5).Value = 1.028440952735571
$ws.Cells.Item(2, 6).Value = 1.044338717707224
$ws.Cells.Item(2, 9).Value = 1.026880918861363
$ws.Cells.Item(2, 10).Value = 1.033672315636142
$ws.Cells.Item(2, 11).Value = 1.039109096347769
$ws.Cells.Item(2, 12).Value = 1.031257761356785
$ws.Cells.Item(2, 13).Value = 1.047110047383149
$ws.Cells.Item(2, 14).Value = 1.035140248888217
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.0297850521199
$ws.Cells.Item(3, 4).Value = 1.037454630074341
$ws.Cells.Item(3, 5).Value = 1.029525195328523
$ws.Cells.Item(3, 6).Value = 1.045554888254531
$ws.Cells.Item(3, 9).Value = 1.026858063055652
$ws.Cells.Item(3, 10).Value = 1.034575943291568
$ws.Cells.Item(3, 11).Value = 1.04005742430864
$ws.Cells.Item(3, 12).Value = 1.03214919475778
$ws.Cells.Item(3, 13).Value = 1.048136379181452
$ws.Cells.Item(3, 14).Value = 1.036045159798562
$ws.Cells.Item(4, 2).Value = 1.019999999999999
$ws.Cells.Item(4, 3).Value = 1.030603310983031
$ws.Cells.Item(4, 4).Value = 1.038192119407687
$ws.Cells.Item(4, 5).Value = 1.030227048339168
$ws.Cells.Item(4, 6).Value = 1.046341380657557
$ws.Cells.Item(4, 9).Value = 1.026840849476427
$ws.Cells.Item(4, 10).Value = 1.035160173170024
$ws.Cells.Item(4, 11).Value = 1.040670553071344
$ws.Cells.Item(4, 12).Value = 1.032725700443172
$ws.Cells.Item(4, 13).Value = 1.048799471330795
$ws.Cells.Item(4, 14).Value = 1.036630219350456
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030947304429489
$ws.Cells.Item(5, 4).Value = 1.038502178963315
$ws.Cells.Item(5, 5).Value = 1.030522174515524
$ws.Cells.Item(5, 6).Value = 1.046671915195816
$ws.Cells.Item(5, 9).Value = 1.026833031603041
$ws.Cells.Item(5, 10).Value = 1.03540567030259
$ws.Cells.Item(5, 11).Value = 1.040928193162715
$ws.Cells.Item(5, 12).Value = 1.03296799001885
$ws.Cells.Item(5, 13).Value = 1.049077993604085
$ws.Cells.Item(5, 14).Value = 1.036876065117111
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03100506236025
$ws.Cells.Item(6, 4).Value = 1.038554240481133
$ws.Cells.Item(6, 5).Value = 1.030571731472968
$ws.Cells.Item(6, 6).Value = 1.046727407166121
$ws.Cells.Item(6, 9).Value = 1.026831684834707
$ws.Cells.Item(6, 10).Value = 1.0354468837546
$ws.Cells.Item(6, 11).Value = 1.040971445108456
$ws.Cells.Item(6, 12).Value = 1.033008667232135
$ws.Cells.Item(6, 13).Value = 1.049124744573417
$ws.Cells.Item(6, 14).Value = 1.036917337096949
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030607907447501
$ws.Cells.Item(7, 4).Value = 1.038196262362247
$ws.Cells.Item(7, 5).Value = 1.030230991564168
$ws.Cells.Item(7, 6).Value = 1.046345797694762
$ws.Cells.Item(7, 9).Value = 1.02684074729882
$ws.Cells.Item(7, 10).Value = 1.035163453957646
$ws.Cells.Item(7, 11).Value = 1.040673996138002
$ws.Cells.Item(7, 12).Value = 1.032728938215438
$ws.Cells.Item(7, 13).Value = 1.048803193907591
$ws.Cells.Item(7, 14).Value = 1.036633504797174
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028947836582486
$ws.Cells.Item(8, 4).Value = 1.036700130437794
$ws.Cells.Item(8, 5).Value = 1.028807321266975
$ws.Cells.Item(8, 6).Value = 1.044749822379461
$ws.Cells.Item(8, 9).Value = 1.026873696079159
$ws.Cells.Item(8, 10).Value = 1.033977800478603
$ws.Cells.Item(8, 11).Value = 1.03942969302709
$ws.Cells.Item(8, 12).Value = 1.031559090587125
$ws.Cells.Item(8, 13).Value = 1.047457111335194
$ws.Cells.Item(8, 14).Value = 1.035446167554194
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026022315250772
$ws.Cells.Item(9, 4).Value = 1.0340642637899
$ws.Cells.Item(9, 5).Value = 1.026300688319896
$ws.Cells.Item(9, 6).Value = 1.041934000713494
$ws.Cells.Item(9, 9).Value = 1.026913227211592
$ws.Cells.Item(9, 10).Value = 1.031884813463793
$ws.Cells.Item(9, 11).Value = 1.037233169145287
$ws.Cells.Item(9, 12).Value = 1.029495234564846
$ws.Cells.Item(9, 13).Value = 1.045077341756215
$ws.Cells.Item(9, 14).Value = 1.033350208257754
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024071639419862
$ws.Cells.Item(10, 4).Value = 1.032307276651265
$ws.Cells.Item(10, 5).Value = 1.024630903393192
$ws.Cells.Item(10, 6).Value = 1.040054348871996
$ws.Cells.Item(10, 9).Value = 1.026927174940004
$ws.Cells.Item(10, 10).Value = 1.030486911090918
$ws.Cells.Item(10, 11).Value = 1.035766120982132
$ws.Cells.Item(10, 12).Value = 1.028117620342596
$ws.Cells.Item(10, 13).Value = 1.04348552400367
$ws.Cells.Item(10, 14).Value = 1.031950320703168
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023226863230478
$ws.Cells.Item(11, 4).Value = 1.031546523065056
$ws.Cells.Item(11, 5).Value = 1.023908158998578
$ws.Cells.Item(11, 6).Value = 1.039239842245172
$ws.Cells.Item(11, 9).Value = 1.026930282394829
$ws.Cells.Item(11, 10).Value = 1.029880973890084
$ws.Cells.Item(11, 11).Value = 1.035130215300255
$ws.Cells.Item(11, 12).Value = 1.027520675570742
$ws.Cells.Item(11, 13).Value = 1.042794976574207
$ws.Cells.Item(11, 14).Value = 1.031343523002007
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022913054578384
$ws.Cells.Item(12, 4).Value = 1.031263948526401
$ws.Cells.Item(12, 5).Value = 1.023639740472647
$ws.Cells.Item(12, 6).Value = 1.038937205597708
$ws.Cells.Item(12, 9).Value = 1.026930996945894
$ws.Cells.Item(12, 10).Value = 1.029655804548066
$ws.Cells.Item(12, 11).Value = 1.034893910228458
$ws.Cells.Item(12, 12).Value = 1.02729887774325
$ws.Cells.Item(12, 13).Value = 1.042538282840253
$ws.Cells.Item(12, 14).Value = 1.031118033893695
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022980368605867
$ws.Cells.Item(13, 4).Value = 1.031324561607772
$ws.Cells.Item(13, 5).Value = 1.0236973153327
$ws.Cells.Item(13, 6).Value = 1.039002126405083
$ws.Cells.Item(13, 9).Value = 1.026930863556654
$ws.Cells.Item(13, 10).Value = 1.029704108566023
$ws.Cells.Item(13, 11).Value = 1.034944603082429
$ws.Cells.Item(13, 12).Value = 1.027346457139954
$ws.Cells.Item(13, 13).Value = 1.042593353304079
$ws.Cells.Item(13, 14).Value = 1.031166406508899
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023200924158653
$ws.Cells.Item(14, 4).Value = 1.031523165310434
$ws.Cells.Item(14, 5).Value = 1.023885970614374
$ws.Cells.Item(14, 6).Value = 1.039214828104966
$ws.Cells.Item(14, 9).Value = 1.026930350424903
$ws.Cells.Item(14, 10).Value = 1.029862363316487
$ws.Cells.Item(14, 11).Value = 1.035110684326917
$ws.Cells.Item(14, 12).Value = 1.02750234304769
$ws.Cells.Item(14, 13).Value = 1.042773762149027
$ws.Cells.Item(14, 14).Value = 1.031324885999261
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023336812819051
$ws.Cells.Item(15, 4).Value = 1.031645531978269
$ws.Cells.Item(15, 5).Value = 1.024002212732351
$ws.Cells.Item(15, 6).Value = 1.039345868317241
$ws.Cells.Item(15, 9).Value = 1.026929976027344
$ws.Cells.Item(15, 10).Value = 1.029959856349038
$ws.Cells.Item(15, 11).Value = 1.035212998979907
$ws.Cells.Item(15, 12).Value = 1.027598380713479
$ws.Cells.Item(15, 13).Value = 1.04288489228739
$ws.Cells.Item(15, 14).Value = 1.031422517483101
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024127701632225
$ws.Cells.Item(16, 4).Value = 1.032357765911716
$ws.Cells.Item(16, 5).Value = 1.02467887538898
$ws.Cells.Item(16, 6).Value = 1.040108392100833
$ws.Cells.Item(16, 9).Value = 1.026926907003469
$ws.Cells.Item(16, 10).Value = 1.030527111597621
$ws.Cells.Item(16, 11).Value = 1.03580830981814
$ws.Cells.Item(16, 12).Value = 1.028157228436325
$ws.Cells.Item(16, 13).Value = 1.043531326279159
$ws.Cells.Item(16, 14).Value = 1.031990578299201
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024623770842143
$ws.Cells.Item(17, 4).Value = 1.032804539119777
$ws.Cells.Item(17, 5).Value = 1.025103402719714
$ws.Cells.Item(17, 6).Value = 1.040586539937908
$ws.Cells.Item(17, 9).Value = 1.026924197477057
$ws.Cells.Item(17, 10).Value = 1.030882764290357
$ws.Cells.Item(17, 11).Value = 1.03618155363576
$ws.Cells.Item(17, 12).Value = 1.028507662974947
$ws.Cells.Item(17, 13).Value = 1.043936473616136
$ws.Cells.Item(17, 14).Value = 1.032346736059555
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024913108181628
$ws.Cells.Item(18, 4).Value = 1.033065137723501
$ws.Cells.Item(18, 5).Value = 1.025351049797073
$ws.Cells.Item(18, 6).Value = 1.040865377229289
$ws.Cells.Item(18, 9).Value = 1.026922334107323
$ws.Cells.Item(18, 10).Value = 1.031090149216139
$ws.Cells.Item(18, 11).Value = 1.036399196469465
$ws.Cells.Item(18, 12).Value = 1.028712024252956
$ws.Cells.Item(18, 13).Value = 1.044172665813425
$ws.Cells.Item(18, 14).Value = 1.032554415495719
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025011762888286
$ws.Cells.Item(19, 4).Value = 1.033153995735219
$ws.Cells.Item(19, 5).Value = 1.025435495759624
$ws.Cells.Item(19, 6).Value = 1.04096044371077
$ws.Cells.Item(19, 9).Value = 1.026921650728612
$ws.Cells.Item(19, 10).Value = 1.031160851774427
$ws.Cells.Item(19, 11).Value = 1.036473396248935
$ws.Cells.Item(19, 12).Value = 1.02878169921563
$ws.Cells.Item(19, 13).Value = 1.044253180361224
$ws.Cells.Item(19, 14).Value = 1.03262521845975
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024570548516945
$ws.Cells.Item(20, 4).Value = 1.032756604250403
$ws.Cells.Item(20, 5).Value = 1.025057852126775
$ws.Cells.Item(20, 6).Value = 1.040535245239177
$ws.Cells.Item(20, 9).Value = 1.026924517444623
$ws.Cells.Item(20, 10).Value = 1.030844612473223
$ws.Cells.Item(20, 11).Value = 1.036141514735822
$ws.Cells.Item(20, 12).Value = 1.028470068945195
$ws.Cells.Item(20, 13).Value = 1.043893017915089
$ws.Cells.Item(20, 14).Value = 1.032308530062465
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.023135976654106
$ws.Cells.Item(21, 4).Value = 1.031464681395861
$ws.Cells.Item(21, 5).Value = 1.023830415219821
$ws.Cells.Item(21, 6).Value = 1.039152195326529
$ws.Cells.Item(21, 9).Value = 1.02693051366078
$ws.Cells.Item(21, 10).Value = 1.029815763936115
$ws.Cells.Item(21, 11).Value = 1.035061780351634
$ws.Cells.Item(21, 12).Value = 1.027456440366522
$ws.Cells.Item(21, 13).Value = 1.042720641598288
$ws.Cells.Item(21, 14).Value = 1.031278220442423
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022233880444361
$ws.Cells.Item(22, 4).Value = 1.030652415924478
$ws.Cells.Item(22, 5).Value = 1.023058912851418
$ws.Cells.Item(22, 6).Value = 1.038282080214551
$ws.Cells.Item(22, 9).Value = 1.02693174016323
$ws.Cells.Item(22, 10).Value = 1.029168322311654
$ws.Cells.Item(22, 11).Value = 1.034382321088581
$ws.Cells.Item(22, 12).Value = 1.026818749450045
$ws.Cells.Item(22, 13).Value = 1.041982400851842
$ws.Cells.Item(22, 14).Value = 1.030629859376593
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022712111476121
$ws.Cells.Item(23, 4).Value = 1.031083012147884
$ws.Cells.Item(23, 5).Value = 1.02346787905059
$ws.Cells.Item(23, 6).Value = 1.038743396156531
$ws.Cells.Item(23, 9).Value = 1.026931330805649
$ws.Cells.Item(23, 10).Value = 1.029511597369375
$ws.Cells.Item(23, 11).Value = 1.03474257151792
$ws.Cells.Item(23, 12).Value = 1.027156838236897
$ws.Cells.Item(23, 13).Value = 1.042373862986745
$ws.Cells.Item(23, 14).Value = 1.03097362192427
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02459459743006
$ws.Cells.Item(24, 4).Value = 1.032778263944243
$ws.Cells.Item(24, 5).Value = 1.02507843439421
$ws.Cells.Item(24, 6).Value = 1.040558423285645
$ws.Cells.Item(24, 9).Value = 1.026924373739448
$ws.Cells.Item(24, 10).Value = 1.030861851826503
$ws.Cells.Item(24, 11).Value = 1.036159606788946
$ws.Cells.Item(24, 12).Value = 1.028487056196542
$ws.Cells.Item(24, 13).Value = 1.043912654057195
$ws.Cells.Item(24, 14).Value = 1.032325793897604
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026778679444768
$ws.Cells.Item(25, 4).Value = 1.034745646895736
$ws.Cells.Item(25, 5).Value = 1.026948477443206
$ws.Cells.Item(25, 6).Value = 1.042662381955979
$ws.Cells.Item(25, 9).Value = 1.026905197884748
$ws.Cells.Item(25, 10).Value = 1.032426349183352
$ws.Cells.Item(25, 11).Value = 1.037801493402287
$ws.Cells.Item(25, 12).Value = 1.030029086875936
$ws.Cells.Item(25, 13).Value = 1.04569349975558
$ws.Cells.Item(25, 14).Value = 1.03389251302015

Write-Output "Applied 264 cell updates"
